$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: Tournevis cruciforme - Entrée
$ws.Cells.Item(24, 1).Value = "2025-05-23 11:20:31"
$ws.Cells.Item(24, 2).Value = "Tournevis cruciforme"
$ws.Cells.Item(24, 3).Value = "Entrée"
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = 29
$ws.Cells.Item(24, 6).Value = 31

# Row 25: Tournevis cruciforme - Sortie
$ws.Cells.Item(25, 1).Value = "2025-05-23 11:21:07"
$ws.Cells.Item(25, 2).Value = "Tournevis cruciforme"
$ws.Cells.Item(25, 3).Value = "Sortie"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 31
$ws.Cells.Item(25, 6).Value = 30
